$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = 4624.449081662315
$ws.Range("D3").Value  = 4624.44908166231
$ws.Range("D4").Value  = 4624.44908166231
$ws.Range("D6").Value  = 13212.71166189251
$ws.Range("D7").Value  = 13212.71166189251
$ws.Range("D9").Value  = 742.3338266865285
$ws.Range("D10").Value = 742.3338266865279
$ws.Range("D11").Value = 10.7957039567556
$ws.Range("D12").Value = 612.7853792054607
$ws.Range("D13").Value = 10.7957039567556
$ws.Range("D17").Value = 8636.563165404481
$ws.Range("D19").Value = 9293.529036356393
$ws.Range("D20").Value = 9293.529036356393
$ws.Range("D21").Value = 129.5484474810672
$ws.Range("D24").Value = 208742.9955873798
$ws.Range("D25").Value = 208742.9955873796
$ws.Range("D28").Value = 4174.859911747636
$ws.Range("D29").Value = 4174.859911747636
$ws.Range("D30").Value = 208742.9955873796
$ws.Range("D35").Value = 20217.59999999929
$ws.Range("D36").Value = 20217.59999999929
$ws.Range("D38").Value = -109.9332197290953
$ws.Range("D39").Value = -109.9332197290952
$ws.Range("D41").Value = 109.9332197290952
$ws.Range("D42").Value = 2198.6643945819
$ws.Range("D43").Value = 2198.6643945819
